# Update column G ("K") values for rows 2-19 in Sheet1, replacing the
# previous Strike# derived values with the regenerated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 3
    17 = 1
    18 = 1
    19 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
